$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2783.1428
$ws.Range("J32").Value = 2579.5
$ws.Range("L32").Value = 2579.5
$ws.Range("N32").Value = -3231.5
$ws.Range("H33").Value = 591.5
$ws.Range("I33").Value = 512.25
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 512.25
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -283.25
$ws.Range("N33").Value = -1208
$ws.Range("H55").Value = 215.88235
$ws.Range("I55").Value = 210.30302
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 210.30302
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = 3.696979999999996
$ws.Range("N55").Value = -828
$ws.Range("H76").Value = 5099.2
$ws.Range("J76").Value = 5098.3335
$ws.Range("L76").Value = 5098.3335
$ws.Range("N76").Value = -5728.3335
$ws.Range("H79").Value = 5099.2
$ws.Range("J79").Value = 5098.3335
$ws.Range("L79").Value = 5098.3335
$ws.Range("N79").Value = -7282.3335
$ws.Range("H80").Value = 322.42856
$ws.Range("I80").Value = 322.42856
$ws.Range("K80").Value = 967.28568
$ws.Range("M80").Value = 30.71432000000004
$ws.Range("H83").Value = 322.42856
$ws.Range("I83").Value = 322.42856
$ws.Range("K83").Value = 2901.85704
$ws.Range("M83").Value = 2090.14296
$ws.Range("H98").Value = 2450.625
$ws.Range("I98").Value = 1900.7858
$ws.Range("K98").Value = 1900.7858
$ws.Range("M98").Value = -402.7858000000001
$ws.Range("H112").Value = 1297.5
$ws.Range("J112").Value = 1297.7567
$ws.Range("L112").Value = 3893.2701
$ws.Range("N112").Value = -6109.2701
$ws.Range("H113").Value = 203874.8
$ws.Range("I113").Value = 501002.5
$ws.Range("J113").Value = 5789.6665
$ws.Range("K113").Value = 501002.5
$ws.Range("L113").Value = 5789.6665
$ws.Range("M113").Value = -497748.5
$ws.Range("N113").Value = -12297.6665
$ws.Range("H122").Value = 2450.625
$ws.Range("I122").Value = 1900.7858
$ws.Range("K122").Value = 5702.357400000001
$ws.Range("M122").Value = -3252.357400000001
$ws.Range("H132").Value = 1780.9762
$ws.Range("J132").Value = 4897.75
$ws.Range("L132").Value = 14693.25
$ws.Range("N132").Value = -19753.25
$ws.Range("H138").Value = 2192.054
$ws.Range("I138").Value = 1691.96
$ws.Range("J138").Value = 3233.9167
$ws.Range("K138").Value = 5075.88
$ws.Range("L138").Value = 9701.750100000001
$ws.Range("M138").Value = 64.11999999999989
$ws.Range("N138").Value = -19981.7501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2465.6
$ws.Range("I102").Value = 1182.5714
$ws.Range("K102").Value = 1182.5714
$ws.Range("M102").Value = 439.4286
$ws.Range("H132").Value = 1519.8975
$ws.Range("I132").Value = 1510.2572
$ws.Range("K132").Value = 4530.7716
$ws.Range("M132").Value = -2000.7716

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5850.4
$ws.Range("I20").Value = 5081.647
$ws.Range("K20").Value = 5081.647
$ws.Range("M20").Value = -4834.647
$ws.Range("H35").Value = 9000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 1539.2858
$ws.Range("I86").Value = 1432.9131
$ws.Range("K86").Value = 1432.9131
$ws.Range("M86").Value = -309.9131
$ws.Range("H89").Value = 1539.2858
$ws.Range("I89").Value = 1432.9131
$ws.Range("K89").Value = 7164.5655
$ws.Range("M89").Value = -1548.5655
$ws.Range("H94").Value = 1544.3954
$ws.Range("J94").Value = 1882.5714
$ws.Range("L94").Value = 1882.5714
$ws.Range("N94").Value = -2784.5714
$ws.Range("H134").Value = 840.8889
$ws.Range("I134").Value = 758.5625
$ws.Range("K134").Value = 2275.6875
$ws.Range("M134").Value = 259.3125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 936.06665
$ws.Range("I16").Value = 1013.5
$ws.Range("J16").Value = 626.3333
$ws.Range("K16").Value = 1013.5
$ws.Range("L16").Value = 626.3333
$ws.Range("M16").Value = -726.5
$ws.Range("N16").Value = -1200.3333
$ws.Range("H62").Value = 5462.125
$ws.Range("J62").Value = 4600
$ws.Range("L62").Value = 4600
$ws.Range("N62").Value = -5848
$ws.Range("H65").Value = 5462.125
$ws.Range("J65").Value = 4600
$ws.Range("L65").Value = 23000
$ws.Range("N65").Value = -29240
$ws.Range("H86").Value = 3785.2727
$ws.Range("I86").Value = 3307.8572
$ws.Range("J86").Value = 4620.75
$ws.Range("K86").Value = 3307.8572
$ws.Range("L86").Value = 4620.75
$ws.Range("M86").Value = -2184.8572
$ws.Range("N86").Value = -6866.75
$ws.Range("H89").Value = 3785.2727
$ws.Range("I89").Value = 3307.8572
$ws.Range("J89").Value = 4620.75
$ws.Range("K89").Value = 16539.286
$ws.Range("L89").Value = 23103.75
$ws.Range("M89").Value = -10923.286
$ws.Range("N89").Value = -34335.75
$ws.Range("H105").Value = 1831.5834
$ws.Range("I105").Value = 1734.4546
$ws.Range("K105").Value = 1734.4546
$ws.Range("M105").Value = 12.54539999999997
$ws.Range("H107").Value = 1444.65
$ws.Range("I107").Value = 1364
$ws.Range("J107").Value = 1553.7646
$ws.Range("K107").Value = 1364
$ws.Range("L107").Value = 1553.7646
$ws.Range("M107").Value = 556
$ws.Range("N107").Value = -5393.7646
$ws.Range("H113").Value = 936.06665
$ws.Range("I113").Value = 1013.5
$ws.Range("J113").Value = 626.3333
$ws.Range("K113").Value = 1013.5
$ws.Range("L113").Value = 626.3333
$ws.Range("M113").Value = 1156.5
$ws.Range("N113").Value = -4966.3333
$ws.Range("H132").Value = 2260.8147
$ws.Range("I132").Value = 2205.8096
$ws.Range("J132").Value = 2453.3333
$ws.Range("K132").Value = 6617.4288
$ws.Range("L132").Value = 7359.999899999999
$ws.Range("M132").Value = -4087.4288
$ws.Range("N132").Value = -12419.9999
$ws.Range("H134").Value = 4010.606
$ws.Range("I134").Value = 3960.5518
$ws.Range("K134").Value = 11881.6554
$ws.Range("M134").Value = -9346.6554

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 363.57144
$ws.Range("I122").Value = 286.2857
$ws.Range("J122").Value = 402.2143
$ws.Range("K122").Value = 2576.5713
$ws.Range("L122").Value = 3619.9287
$ws.Range("M122").Value = -126.5713000000001
$ws.Range("N122").Value = -8519.9287
$ws.Range("H131").Value = 60130.766
$ws.Range("I131").Value = 77824.84
$ws.Range("J131").Value = 2625
$ws.Range("K131").Value = 233474.52
$ws.Range("L131").Value = 7875
$ws.Range("M131").Value = -228434.52
$ws.Range("N131").Value = -17955
$ws.Range("H132").Value = 1756.091
$ws.Range("I132").Value = 1157
$ws.Range("J132").Value = 2035.6666
$ws.Range("K132").Value = 10413
$ws.Range("L132").Value = 18320.9994
$ws.Range("M132").Value = -7883
$ws.Range("N132").Value = -23380.9994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H70").Value = 5169.5
$ws.Range("I70").Value = 5111.625
$ws.Range("K70").Value = 5111.625
$ws.Range("M70").Value = -4841.625
$ws.Range("H73").Value = 5169.5
$ws.Range("I73").Value = 5111.625
$ws.Range("K73").Value = 5111.625
$ws.Range("M73").Value = -4175.625
$ws.Range("H113").Value = 2623.75
$ws.Range("J113").Value = 2865
$ws.Range("L113").Value = 2865
$ws.Range("N113").Value = -7205
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H126").Value = 3546.4614
$ws.Range("I126").Value = 3072.2856
$ws.Range("K126").Value = 9216.856800000001
$ws.Range("M126").Value = -6746.856800000001
$ws.Range("H132").Value = 3445.9697
$ws.Range("I132").Value = 3812.3333
$ws.Range("K132").Value = 11436.9999
$ws.Range("M132").Value = -8906.999899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 38161.582
$ws.Range("J46").Value = 4074.1428
$ws.Range("L46").Value = 4074.1428
$ws.Range("N46").Value = -4450.1428
$ws.Range("H55").Value = 613.44446
$ws.Range("I55").Value = 620
$ws.Range("K55").Value = 620
$ws.Range("M55").Value = -447
$ws.Range("H122").Value = 6811.577
$ws.Range("J122").Value = 3732.375
$ws.Range("L122").Value = 11197.125
$ws.Range("N122").Value = -16097.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 12500
$ws.Range("J25").Value = 12500
$ws.Range("L25").Value = 12500
$ws.Range("N25").Value = -13086
$ws.Range("H62").Value = 20374.438
$ws.Range("I62").Value = 14539.2
$ws.Range("J62").Value = 23026.818
$ws.Range("K62").Value = 14539.2
$ws.Range("L62").Value = 23026.818
$ws.Range("M62").Value = -13915.2
$ws.Range("N62").Value = -24274.818
$ws.Range("H65").Value = 20374.438
$ws.Range("I65").Value = 14539.2
$ws.Range("J65").Value = 23026.818
$ws.Range("K65").Value = 72696
$ws.Range("L65").Value = 115134.09
$ws.Range("M65").Value = -69576
$ws.Range("N65").Value = -121374.09
$ws.Range("H86").Value = 16330
$ws.Range("J86").Value = 16330
$ws.Range("L86").Value = 16330
$ws.Range("N86").Value = -18576
$ws.Range("H89").Value = 16330
$ws.Range("J89").Value = 16330
$ws.Range("L89").Value = 81650
$ws.Range("N89").Value = -92882
$ws.Range("H137").Value = 30049358
$ws.Range("J137").Value = 30049358
$ws.Range("L137").Value = 30049358
$ws.Range("N137").Value = -30059558
